$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" (sheet 1) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Remove the "VERA CABRERA JORGE ENRIQUE" row (row 59). Everything below
# shifts up one row: the old row 60 (VIVANCO MALDONADO SILVANA MARILY)
# becomes the new row 59, and the old totals row 61 becomes the new row 60.
$ws1.Rows.Item(59).Delete()

# The totals row (now row 60) counted "out of 59" advisors; after removing
# one advisor row it must read "out of 58". The counts themselves are
# unchanged, only the denominator text.
$ws1.Range("C60").Value = "0 de 58"
$ws1.Range("D60").Value = "1 de 58"
$ws1.Range("E60").Value = "1 de 58"
$ws1.Range("F60").Value = "0 de 58"
$ws1.Range("G60").Value = "0 de 58"
$ws1.Range("H60").Value = "1 de 58"
$ws1.Range("I60").Value = "1 de 58"
$ws1.Range("J60").Value = "0 de 58"
$ws1.Range("K60").Value = "0 de 58"
$ws1.Range("L60").Value = "0 de 58"
$ws1.Range("M60").Value = "4 de 58"
$ws1.Range("N60").Value = "0 de 58"
$ws1.Range("O60").Value = "2 de 58"
$ws1.Range("P60").Value = "0 de 58"
$ws1.Range("Q60").Value = "0 de 58"
$ws1.Range("R60").Value = "0 de 58"

# --- Sheet "VENTA MENSUAL" (sheet 2) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Same row removal as above: drop the "VERA CABRERA JORGE ENRIQUE" row.
$ws2.Rows.Item(59).Delete()

# Recompute the grand-total row (now row 60) now that one advisor's
# figures have been removed from the sums.
$ws2.Range("C60").Value = 53770.63
$ws2.Range("D60").Value = 53165.42
$ws2.Range("E60").Value = 57170.26
$ws2.Range("F60").Value = 11393.66
$ws2.Range("G60").Value = 47000
